$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays formatted as Text so numeric-looking
# strings like "1.002" or "27.445.73" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Simple D/E column updates
$ws.Range("D2").Value = "27.445.73"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.736.30"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "322.91"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "0.4542"
$ws.Range("E7").Value = "  +7.58%  "
$ws.Range("D8").Value = "0.3526"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").Value = "0.07399"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "41.34"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").Value = "1.070"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "20.33"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "5.894"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D17").Value = "91.17"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "0.06339"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").Value = "16.52"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").Value = "5.709"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "27.500.26"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "161.85"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "19.95"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "1.936.89"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D31").Value = "1.043"
$ws.Range("E31").Value = "  -5.71%  "
$ws.Range("D32").Value = "0.09072"
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("D33").Value = "3.648"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "5.364"
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("D35").Value = "0.02265"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "11.57"
$ws.Range("E36").Value = "  -5.22%  "
$ws.Range("D37").Value = "0.05935"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "0.2052"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "0.6216"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "4.869"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "1.185"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "1.375"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").Value = "7.668"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").Value = "13.11"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "3.701"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "0.5770"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "122.02"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "1.926"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").Value = "0.06842"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").Value = "1.109"
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("D51").Value = "71.00"
$ws.Range("E51").Value = "  -2.71%  "

# Row swap: 15 (Chainlink) <-> 16 (WrappedEther)
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.742.83"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.034"
$ws.Range("E16").Value = "  -2.27%  "

# Row swap: 29 (LidoDAOToken) <-> 30 (BitcoinCash)
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "124.44"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "2.035"
$ws.Range("E30").Value = "  -3.97%  "
